$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.257.43"

$ws.Range("D3").Value = "1.868.33"
$ws.Range("E3").Value = "  +1.15%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'337.96"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").Value = "'0.4704"
$ws.Range("E7").Value = "  +1.47%  "

$ws.Range("D8").Value = "'0.3922"
$ws.Range("E8").Value = "  +1.75%  "

$ws.Range("D9").Value = "'47.21"
$ws.Range("E9").Value = "  +2.35%  "

$ws.Range("D10").Value = "'0.07994"
$ws.Range("E10").Value = "  +0.93%  "

$ws.Range("D11").Value = "'1.004"
$ws.Range("E11").Value = "  +0.95%  "

$ws.Range("D12").Value = "'21.73"
$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("D13").Value = "1.876.29"
$ws.Range("E13").Value = "  +1.58%  "

$ws.Range("E14").Value = "  +1.15%  "

$ws.Range("D15").Value = "'7.263"
$ws.Range("E15").Value = "  +2.17%  "

$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'91.08"
$ws.Range("E17").Value = "  +2.48%  "

$ws.Range("D18").Value = "'0.00001042"
$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("D19").Value = "'0.06582"
$ws.Range("E19").Value = "  -0.80%  "

$ws.Range("D20").Value = "'17.61"
$ws.Range("E20").Value = "  +3.11%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("D22").Value = "28.272.77"
$ws.Range("E22").Value = "  +2.63%  "

$ws.Range("D23").Value = "'5.440"
$ws.Range("E23").Value = "  +1.06%  "

$ws.Range("D24").Value = "'11.04"
$ws.Range("E24").Value = "  +1.09%  "

$ws.Range("E25").Value = "  -0.31%  "

$ws.Range("D26").Value = "2.087.36"
$ws.Range("E26").Value = "  +1.09%  "

$ws.Range("D27").Value = "'159.72"
$ws.Range("E27").Value = "  +1.00%  "

$ws.Range("D28").Value = "'19.79"
$ws.Range("E28").Value = "  +1.35%  "

$ws.Range("D29").Value = "'2.141"
$ws.Range("E29").Value = "  +1.82%  "

$ws.Range("D30").Value = "'5.485"
$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("D31").Value = "'119.95"
$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("D32").Value = "'0.9751"
$ws.Range("E32").Value = "  -0.11%  "

$ws.Range("D33").Value = "'0.09498"
$ws.Range("E33").Value = "  +0.96%  "

$ws.Range("D34").Value = "'3.591"
$ws.Range("E34").Value = "  +0.28%  "

$ws.Range("D35").Value = "'1.373"
$ws.Range("E35").Value = "  +1.87%  "

$ws.Range("D36").Value = "'5.343"
$ws.Range("E36").Value = "  +1.09%  "

$ws.Range("D37").Value = "'0.02273"
$ws.Range("E37").Value = "  +2.22%  "

$ws.Range("D38").Value = "'0.06083"
$ws.Range("E38").Value = "  +1.13%  "

$ws.Range("D39").Value = "'8.419"
$ws.Range("E39").Value = "  +1.52%  "

$ws.Range("D40").Value = "'1.176"
$ws.Range("E40").Value = "  -0.36%  "

$ws.Range("D41").Value = "'0.5961"
$ws.Range("E41").Value = "  +1.18%  "

$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("D43").Value = "'0.1876"
$ws.Range("E43").Value = "  +0.79%  "

$ws.Range("D44").Value = "'10.36"
$ws.Range("E44").Value = "  +0.52%  "

$ws.Range("D45").Value = "'1.301"
$ws.Range("E45").Value = "  +4.57%  "

$ws.Range("D46").Value = "'0.5602"
$ws.Range("E46").Value = "  +0.37%  "

$ws.Range("D47").Value = "'12.17"
$ws.Range("E47").Value = "  +0.29%  "

$ws.Range("D48").Value = "'1.963"
$ws.Range("E48").Value = "  +3.39%  "

$ws.Range("D49").Value = "'0.06897"
$ws.Range("E49").Value = "  +3.27%  "

$ws.Range("D50").Value = "'110.81"
$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D51").Value = "'1.997"
$ws.Range("E51").Value = "  +12.18%  "
